$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the text of cell A6: fix wording back to "Manejoa los errores"
$ws.Range("A6").Value = "Cumple con Clean Architecture (Plural,Metodos CRUD, Manejoa los errores)"

# Move the active selection from A6 to C6
$ws.Range("C6").Select()
